# Refresh the scraped crypto snapshot (price + 1h volume delta columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume cells hold numeric-looking text (e.g. "35.70", "72.30").
# A leading apostrophe is the standard Excel way to force text entry so
# the engine does not reinterpret the literal as a Number and silently
# normalise it (dropping trailing zeros, re-parsing "." as thousands, etc).

$ws.Range("D2").Value = '''48.125.99'
$ws.Range("E2").Value = '''  +1.81%  '
$ws.Range("D3").Value = '''2.524.81'
$ws.Range("E3").Value = '''  +0.99%  '
$ws.Range("E4").Value = '''  -0.05%  '
$ws.Range("D5").Value = '''323.85'
$ws.Range("E5").Value = '''  +0.10%  '
$ws.Range("D6").Value = '''108.96'
$ws.Range("E6").Value = '''  -0.02%  '
$ws.Range("E8").Value = '''  -0.01%  '
$ws.Range("D9").Value = '''0.558'
$ws.Range("E9").Value = '''  +4.15%  '
$ws.Range("D10").Value = '''40.58'
$ws.Range("E10").Value = '''  +3.59%  '
$ws.Range("D11").Value = '''20.40'
$ws.Range("E11").Value = '''  +10.96%  '
$ws.Range("E12").Value = '''  +1.14%  '
$ws.Range("E13").Value = '''  +1.11%  '
$ws.Range("E14").Value = '''  +1.07%  '
$ws.Range("D15").Value = '''2.919.06'
$ws.Range("E15").Value = '''  +0.92%  '
$ws.Range("D16").Value = '''2.525.43'
$ws.Range("E16").Value = '''  +1.03%  '
$ws.Range("D17").Value = '''0.860'
$ws.Range("E17").Value = '''  +0.72%  '
$ws.Range("D18").Value = '''47.976.73'
$ws.Range("E18").Value = '''  +1.63%  '
$ws.Range("D19").Value = '''13.25'
$ws.Range("E19").Value = '''  +3.20%  '
$ws.Range("D20").Value = '''6.63'
$ws.Range("D21").Value = '''0.0₃0952'
$ws.Range("E21").Value = '''  +1.17%  '
$ws.Range("E22").Value = '''  -0.48%  '
$ws.Range("D23").Value = '''72.30'
$ws.Range("E23").Value = '''  +2.17%  '
$ws.Range("D24").Value = '''269.72'
$ws.Range("E24").Value = '''  +8.82%  '
$ws.Range("E25").Value = '''  -0.74%  '
$ws.Range("D26").Value = '''26.21'
$ws.Range("E26").Value = '''  +0.54%  '
$ws.Range("E27").Value = '''  -0.28%  '
$ws.Range("D28").Value = '''10.16'
$ws.Range("E28").Value = '''  +0.93%  '
$ws.Range("E29").Value = '''  +4.87%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '''35.70'
$ws.Range("E30").Value = '''  +1.35%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '''2.10'
$ws.Range("E31").Value = '''  -8.38%  '
$ws.Range("D32").Value = '''49.79'
$ws.Range("E32").Value = '''  -0.01%  '
$ws.Range("D33").Value = '''19.97'
$ws.Range("E33").Value = '''  -0.15%  '
$ws.Range("E34").Value = '''  -0.37%  '
$ws.Range("E35").Value = '''  -0.07%  '
$ws.Range("E36").Value = '''  +0.78%  '
$ws.Range("D37").Value = '''2.00'
$ws.Range("E37").Value = '''  +0.80%  '
$ws.Range("E38").Value = '''  +0.73%  '
$ws.Range("E39").Value = '''  +0.40%  '
$ws.Range("E40").Value = '''  -0.06%  '
$ws.Range("D41").Value = '''22.46'
$ws.Range("E41").Value = '''  +5.59%  '
$ws.Range("D42").Value = '''2.19'
$ws.Range("E42").Value = '''  -2.00%  '
$ws.Range("D43").Value = '''119.04'
$ws.Range("E43").Value = '''  -2.51%  '
$ws.Range("E44").Value = '''  +0.45%  '
$ws.Range("D45").Value = '''2.011.92'
$ws.Range("E45").Value = '''  +1.00%  '
$ws.Range("D46").Value = '''3.13'
$ws.Range("E46").Value = '''  +2.50%  '
$ws.Range("D47").Value = '''2.06'
$ws.Range("E47").Value = '''  -0.03%  '
$ws.Range("D48").Value = '''1.88'
$ws.Range("E48").Value = '''  +5.08%  '
$ws.Range("D49").Value = '''9.15'
$ws.Range("E49").Value = '''  +0.73%  '
$ws.Range("E50").Value = '''  +0.71%  '
$ws.Range("D51").Value = '''79.95'
$ws.Range("E51").Value = '''  +2.51%  '
